$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-03-08 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-03-09 Sunday", 2) | Out-Null
$d.Content.Find.Execute("120÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "116÷8=", 2) | Out-Null
$d.Content.Find.Execute("122÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "852÷8=", 2) | Out-Null
$d.Content.Find.Execute("856÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "526÷3=", 2) | Out-Null
$d.Content.Find.Execute("383÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "971÷6=", 2) | Out-Null
$d.Content.Find.Execute("272÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "134÷9=", 2) | Out-Null
$d.Content.Find.Execute("876÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "346÷3=", 2) | Out-Null
$d.Content.Find.Execute("225÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "464÷8=", 2) | Out-Null
$d.Content.Find.Execute("701÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "509÷7=", 2) | Out-Null
$d.Content.Find.Execute("251÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "404÷8=", 2) | Out-Null
$d.Content.Find.Execute("679÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "915÷9=", 2) | Out-Null
$d.Content.Find.Execute("189÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "605÷2=", 2) | Out-Null
$d.Content.Find.Execute("945÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "601÷8=", 2) | Out-Null
$d.Content.Find.Execute("870÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "165÷2=", 2) | Out-Null
$d.Content.Find.Execute("754÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "691÷9=", 2) | Out-Null
$d.Content.Find.Execute("182÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "475÷9=", 2) | Out-Null
$d.Content.Find.Execute("675÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "629÷6=", 2) | Out-Null
$d.Content.Find.Execute("344÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "180÷6=", 2) | Out-Null
$d.Content.Find.Execute("881÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "596÷5=", 2) | Out-Null
$d.Content.Find.Execute("208÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "986÷8=", 2) | Out-Null
$d.Content.Find.Execute("418÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "911÷5=", 2) | Out-Null
$d.Content.Find.Execute("101÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "952÷2=", 2) | Out-Null
$d.Content.Find.Execute("846÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "653÷5=", 2) | Out-Null
$d.Content.Find.Execute("462÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "104÷3=", 2) | Out-Null
$d.Content.Find.Execute("546÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "294÷8=", 2) | Out-Null
$d.Content.Find.Execute("655÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "155÷2=", 2) | Out-Null
